$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("train")

# New hole_id header and values to replace the old numeric index in column A
$holeIds = @(
    "LBU_02_3",
    "LBU_05_23",
    "MHZ_12_03",
    "LBU_05_17",
    "LBU_05_06",
    "MHZ_08_02",
    "LBU_05_15",
    "LBU_05_29",
    "MHZ_12_01",
    "LBU_05_13",
    "MHZ_08_01",
    "LBU_98_6",
    "LBU_96_4",
    "LBU_05_27",
    "LBU_98_1",
    "LBU_87_1",
    "MHZ_08_05",
    "LBU_05_09",
    "LBU_05_18",
    "LBU_05_28",
    "LBU_05_26",
    "LBU_07_01",
    "LBU_05_19",
    "LBU_02_4",
    "LBU_05_10",
    "LBU_05_25",
    "LBU_05_24",
    "MHZ_08_04",
    "LBU_96_1",
    "LBU_05_14",
    "MHZ_12_04",
    "LBU_05_16",
    "LBU_05_21",
    "LBU_05_02",
    "LBU_96_2",
    "LBU_05_03",
    "MHZ_12_02",
    "LBU_01_2",
    "LBU_05_12",
    "LBU_01_1",
    "LBU_05_11",
    "LBU_05_01",
    "LBU_05_07",
    "LBU_05_05",
    "LBU_05_04",
    "LBU_05_08",
    "LBU_87_6",
    "LBU_07_03"
)

$ws.Range("A1").Value = "hole_id"
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null

for ($i = 0; $i -lt $holeIds.Length; $i++) {
    $row = $i + 2
    $ws.Cells.Item($row, 1).Value = $holeIds[$i]
}
